# Refactored portfolio to use instrument
#
# - Column E ("Category *" / Unlisted) is repurposed in place to become
#   "Instrument" / "Common Stock".
# - Column F ("Sub Category *" / Equity) is no longer needed and is removed.
# - The comment that lives on E1 is kept (anchored to E1) but its text is
#   updated; the comment that lived on F1 is deleted outright.
# - The selection is moved to G10 (matches the saved view state in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the comment anchored at E1 (stays on E1 even though the column's
# meaning changes) to reflect the new "Instrument" column.
$null = $ws.Range("E1").Comment.Text("Author:
-Mandatory
-This is the instrument of the portfolio company being valued")

# Remove the comment that was on F1 ("Sub Category *") since that column is
# going away entirely.
$null = $ws.Range("F1").Comment.Delete()

# Repurpose column E's header and values in place.
$ws.Range("E1").Value = "Instrument"
$ws.Range("E2:E5").Value = "Common Stock"

# Column F ("Sub Category *" / Equity) is no longer used - delete it.
$null = $ws.Range("F1").EntireColumn.Delete()

# Restore the saved selection state from the workbook.
$null = $ws.Range("G10").Select()
